$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new weekly record at row 59 (pushes former rows 59-75 down to 60-76)
$ws.Rows.Item(59).Insert()
$ws.Range("A59").Value = 6
$ws.Range("B59").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = 44782
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100108
$ws.Range("H59").Value = "Tropicales y subtropicales"
$ws.Range("I59").Value = 100108007
$ws.Range("J59").Value = "Coco"
$ws.Range("K59").Value = "Sin especificar"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 150
$ws.Range("N59").Value = 28000
$ws.Range("O59").Value = 29000
$ws.Range("P59").Value = 28500
$ws.Range("Q59").Value = "$/malla 20 unidades"
$ws.Range("R59").Value = "Perú"
$ws.Range("S59").Value = 1425
$ws.Range("T59").Value = 20

# Insert second new weekly record at row 70 (pushes former rows 69-75, already
# shifted once to 70-76, further down to 71-77)
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = 6
$ws.Range("B70").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = 44783
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100108
$ws.Range("H70").Value = "Tropicales y subtropicales"
$ws.Range("I70").Value = 100108007
$ws.Range("J70").Value = "Coco"
$ws.Range("K70").Value = "Sin especificar"
$ws.Range("L70").Value = "Primera"
$ws.Range("M70").Value = 100
$ws.Range("N70").Value = 28000
$ws.Range("O70").Value = 29000
$ws.Range("P70").Value = 28500
$ws.Range("Q70").Value = "$/malla 20 unidades"
$ws.Range("R70").Value = "Perú"
$ws.Range("S70").Value = 1425
$ws.Range("T70").Value = 20
